$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the rule-table header cell F7 from "autogen" to "autogen == $param"
$ws.Range("F7").Value = "autogen == `$param"

# Move the active cell/selection to F7 (matches recorded cursor position in the saved file)
$ws.Activate()
$ws.Range("F7").Select()
